$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the price cells that are plain decimal numbers as Text first,
# so assigning the literal string below keeps the exact digits/trailing
# zeros (matching the raw text feed) instead of Excel coercing them
# into numeric values.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D33", "D38", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.482.69'
$ws.Range("E2").Value = '  +2.84%  '

$ws.Range("D3").Value = '1.603.55'
$ws.Range("E3").Value = '  +2.90%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '213.09'
$ws.Range("E5").Value = '  +1.06%  '

$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +6.57%  '

$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("D8").Value = '27.06'
$ws.Range("E8").Value = '  +10.32%  '

$ws.Range("D9").Value = '43.53'
$ws.Range("E9").Value = '  -4.39%  '

$ws.Range("E10").Value = '  +2.00%  '

$ws.Range("D11").Value = '0.0598'
$ws.Range("E11").Value = '  +2.47%  '

$ws.Range("D12").Value = '0.0910'
$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("D13").Value = '1.832.66'
$ws.Range("E13").Value = '  +3.02%  '

$ws.Range("D14").Value = '1.604.67'
$ws.Range("E14").Value = '  +3.03%  '

$ws.Range("D15").Value = '29.532.96'
$ws.Range("E15").Value = '  +3.11%  '

$ws.Range("E16").Value = '  +4.45%  '

$ws.Range("E17").Value = '  +2.49%  '

$ws.Range("D18").Value = '63.51'
$ws.Range("E18").Value = '  +3.56%  '

$ws.Range("D19").Value = '243.59'
$ws.Range("E19").Value = '  +5.50%  '

$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +3.01%  '

$ws.Range("E21").Value = '  +2.88%  '

$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = '4.02'
$ws.Range("E23").Value = '  +2.60%  '

$ws.Range("D24").Value = '9.16'
$ws.Range("E24").Value = '  +2.05%  '

$ws.Range("D25").Value = '2.07'
$ws.Range("E25").Value = '  -1.98%  '

$ws.Range("D26").Value = '154.42'
$ws.Range("E26").Value = '  +1.83%  '

$ws.Range("E27").Value = '  +3.48%  '

$ws.Range("E28").Value = '  +5.10%  '

$ws.Range("D29").Value = '6.38'
$ws.Range("E29").Value = '  +2.07%  '

$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +2.39%  '

$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("D33").Value = '3.22'
$ws.Range("E33").Value = '  +1.55%  '

$ws.Range("D34").Value = '1.420.90'
$ws.Range("E34").Value = '  +1.94%  '

$ws.Range("E35").Value = '  +3.23%  '

$ws.Range("E36").Value = '  -2.37%  '

$ws.Range("E37").Value = '  +1.85%  '

$ws.Range("D38").Value = '2.81'
$ws.Range("E38").Value = '  +4.38%  '

$ws.Range("E39").Value = '  +1.14%  '

$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("D41").Value = '0.533'
$ws.Range("E41").Value = '  +2.87%  '

$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").Value = '52.70'
$ws.Range("E44").Value = '  +19.94%  '

$ws.Range("D45").Value = '0.792'
$ws.Range("E45").Value = '  +2.00%  '

$ws.Range("D46").Value = '0.0473'
$ws.Range("E46").Value = '  +1.67%  '

$ws.Range("D47").Value = '65.85'
$ws.Range("E47").Value = '  +2.35%  '

$ws.Range("D48").Value = '5.30'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").Value = '1.746.55'
$ws.Range("E49").Value = '  +3.00%  '

$ws.Range("D50").Value = '86.35'
$ws.Range("E50").Value = '  +1.17%  '

$ws.Range("D51").Value = '0.833'
$ws.Range("E51").Value = '  -3.88%  '
